# Update Name of Algo
# Apply updated KNN-imputed values to the result data sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "E3"  = 16.209
    "C7"  = -12.937
    "B10" = 5.922
    "B12" = 5.315
    "C15" = -13.636
    "B18" = 5.112
    "E18" = 16.601
    "E19" = 16.538
    "C20" = -12.183
    "E27" = 16.357
    "C29" = -11.97
    "C30" = -13.177
    "C31" = -13.185
    "B37" = 8.494
    "C40" = -12.782
    "E42" = 16.466
    "E44" = 16.542
    "E47" = 16.399
    "B55" = 5.315
    "E58" = 16.541
    "B68" = 5.431
    "C68" = -11.177
    "E73" = 16.561
    "C76" = -12.519
    "B77" = 6.377
    "B78" = 7.58
    "C87" = -13.071
    "C88" = -13.089
    "E95" = 17.602
    "C96" = -12.57
    "C98" = -13.201
    "C101" = -12.612
    "E101" = 16.655
    "C102" = -13.091
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
